$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 = "Save", copy formatting (style) from the existing G1 header cell
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" data column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
